$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (one month later: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the NEGRO price list (rows 26-32) and BLANCO price list (rows 34-40)
# with the new prices (values read ~1.9699x the previous ones).
$ws.Range("D26").Value = 6797.942
$ws.Range("D27").Value = 8507.816999999999
$ws.Range("D28").Value = 11071.228
$ws.Range("D29").Value = 15643.851
$ws.Range("D30").Value = 23056.99
$ws.Range("D31").Value = 31703.345
$ws.Range("D32").Value = 40757.109

$ws.Range("D34").Value = 6797.942
$ws.Range("D35").Value = 8507.816999999999
$ws.Range("D36").Value = 11071.228
$ws.Range("D37").Value = 15643.851
$ws.Range("D38").Value = 23056.99
$ws.Range("D39").Value = 31703.345
$ws.Range("D40").Value = 40757.109
